$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width changes (cols P=16, S=19, T=20, U=21)
# Note: engine quantizes ColumnWidth to nearest 1/6 character unit internally,
# so we pick the closest settable ColumnWidth to reach each target stored width.
$ws.Columns.Item(16).ColumnWidth = 20.833333333333332
$ws.Columns.Item(19).ColumnWidth = 18.833333333333332
$ws.Columns.Item(20).ColumnWidth = 22.833333333333332
$ws.Columns.Item(21).ColumnWidth = 18.833333333333332

# Cell value changes (error calculation results)
$ws.Range("B2").Value = -105.543252184767
$ws.Range("C2").Value = 0.001522748921522226
$ws.Range("L2").Value = 0.006779213643547647
$ws.Range("M2").Value = 0.1702428792942347
$ws.Range("N2").Value = 0.00004916713429368548
$ws.Range("O2").Value = 0.1702428792942349
$ws.Range("P2").Value = 1.318826771892574
$ws.Range("Q2").Value = 0.5335004453826736
$ws.Range("T2").Value = 0.2257985771763353
$ws.Range("U2").Value = 1.726512574742643
$ws.Range("C3").Value = 0.002107527865453386
$ws.Range("D3").Value = 0.01001464435336744
$ws.Range("E3").Value = 3.988986395410593
$ws.Range("F3").Value = 0.007159443268110569
$ws.Range("G3").Value = 0.2433412122966173
$ws.Range("H3").Value = 0.2087949818774892
$ws.Range("I3").Value = 2.784202105567774
$ws.Range("J3").Value = 0.157989250425516
$ws.Range("K3").Value = 2.214395630898379
$ws.Range("M3").Value = 0.2105665503060151
$ws.Range("O3").Value = 0.2105665503060163
$ws.Range("P3").Value = 0.03525043890623247
$ws.Range("Q3").Value = 5.705288904693922
$ws.Range("R3").Value = 0.3294613250209388
$ws.Range("S3").Value = 1.213265190703177
$ws.Range("T3").Value = 0.01157264535611715
$ws.Range("U3").Value = 2.453371522281945
$ws.Range("P4").Value = 1.318759177757272
$ws.Range("Q4").Value = 0.6996313538309004
$ws.Range("T4").Value = 0.2271470050020687
$ws.Range("U4").Value = 1.004207089032775
$ws.Range("B5").Value = -5.745885478354085
$ws.Range("C5").Value = 0.005078261061152813
$ws.Range("D5").Value = 0.0100333743049465
$ws.Range("E5").Value = 0.4349619753938642
$ws.Range("F5").Value = 0.007142862075553487
$ws.Range("G5").Value = 0.1244576487916852
$ws.Range("H5").Value = 0.0340962603054759
$ws.Range("I5").Value = 0.09998252009169473
$ws.Range("K5").Value = 0.6842968537354914
$ws.Range("L5").Value = 0.007535591938662253
$ws.Range("M5").Value = 0.5107608796364962
$ws.Range("N5").Value = 0.00005465286688276307
$ws.Range("O5").Value = 0.5107608796364961
$ws.Range("P5").Value = 0.005921825770778144
$ws.Range("Q5").Value = 2.813032265449945
$ws.Range("R5").Value = 0.2259939097674143
$ws.Range("S5").Value = 1.381875981856847
$ws.Range("T5").Value = 0.001274691673562648
$ws.Range("U5").Value = 2.97747271498213
$ws.Range("E6").Value = 2.657062243021936
$ws.Range("F6").Value = 0.007138611252906625
$ws.Range("G6").Value = 0.01676710101204121
$ws.Range("H6").Value = 2.732088290507844
$ws.Range("I6").Value = 0.03758627650321703
$ws.Range("P6").Value = 1.313374031600917
$ws.Range("Q6").Value = 0.5419107723895897
$ws.Range("T6").Value = 0.2393293430170869
$ws.Range("U6").Value = 1.28265146190901
$ws.Range("B7").Value = 9.546339218875932
$ws.Range("C7").Value = 0.001978891721470877
$ws.Range("D7").Value = 0.01001478293103228
$ws.Range("E7").Value = 0.2991249825553106
$ws.Range("F7").Value = 0.007148360260622351
$ws.Range("G7").Value = 0.03217013055235829
$ws.Range("H7").Value = 0.13025665020751
$ws.Range("I7").Value = 0.0304533949071753
$ws.Range("J7").Value = 0.09979483482092993
$ws.Range("K7").Value = 0.4047604426046078
$ws.Range("L7").Value = 0.007651493862998872
$ws.Range("M7").Value = 0.1960179186031241
$ws.Range("N7").Value = 0.00005549346075963237
$ws.Range("O7").Value = 0.1960179186031242
$ws.Range("P7").Value = 0.01783305396177997
$ws.Range("Q7").Value = 1.001748480524192
$ws.Range("R7").Value = 0.249544100283467
$ws.Range("S7").Value = 0.5610059207786317
$ws.Range("T7").Value = 0.004458367135931112
$ws.Range("U7").Value = 0.6951876357866389
$ws.Range("B8").Value = -75.86083758120998
$ws.Range("C8").Value = 0.001344655020757129
$ws.Range("F8").Value = 0.007147061784218761
$ws.Range("G8").Value = 0.0136597412716
$ws.Range("J8").Value = 1.909622255286029
$ws.Range("K8").Value = 2.058287584107706
$ws.Range("L8").Value = 0.007004180843521679
$ws.Range("M8").Value = 0.145503521053876
$ws.Range("N8").Value = 0.00005079873835787148
$ws.Range("O8").Value = 0.1455035210538751
$ws.Range("P8").Value = 1.313892579368417
$ws.Range("Q8").Value = 0.4004394967171863
$ws.Range("T8").Value = 0.2374191339979786
$ws.Range("U8").Value = 1.300845587375957
$ws.Range("B9").Value = 12.36075644313317
$ws.Range("C9").Value = 0.001617865900988797
$ws.Range("I9").Value = 0.008632726519021489
$ws.Range("J9").Value = 0.2950203981457464
$ws.Range("K9").Value = 0.2452953080741899
$ws.Range("L9").Value = 0.007672824727450308
$ws.Range("M9").Value = 0.1598112027448657
$ws.Range("N9").Value = 0.00005564816564610286
$ws.Range("O9").Value = 0.1598112027448663
$ws.Range("P9").Value = 0.1438931916272418
$ws.Range("Q9").Value = 0.365815313729743
$ws.Range("T9").Value = 0.002664781889775179
$ws.Range("U9").Value = 0.3003130595442578
$ws.Range("B10").Value = -71.56286999932171
$ws.Range("C10").Value = 0.001234439821687892
$ws.Range("D10").Value = 0.009704301232877284
$ws.Range("E10").Value = 2.952408997433019
$ws.Range("G10").Value = 0.01643770650865622
$ws.Range("H10").Value = 2.730883348649252
$ws.Range("I10").Value = 0.02776156534325191
$ws.Range("J10").Value = 1.978969540598102
$ws.Range("K10").Value = 2.995900406725089
$ws.Range("L10").Value = 0.007036755745037969
$ws.Range("M10").Value = 0.1329589028485945
$ws.Range("N10").Value = 0.0000510349920949077
$ws.Range("O10").Value = 0.1329589028485957
$ws.Range("P10").Value = 1.312935307459668
$ws.Range("Q10").Value = 0.5991393735161512
$ws.Range("T10").Value = 0.2340798179609642
$ws.Range("U10").Value = 1.249551913334694
$ws.Range("C11").Value = 0.001280124441955429
$ws.Range("D11").Value = 0.009527788728580266
$ws.Range("E11").Value = 0.3030891627036431
$ws.Range("H11").Value = 0.4493095712755549
$ws.Range("I11").Value = 0.009845085273641255
$ws.Range("M11").Value = 0.1181142656994033
$ws.Range("O11").Value = 0.1181142656994044
$ws.Range("P11").Value = 0.1434984853044898
$ws.Range("Q11").Value = 0.3364171721037157
$ws.Range("R11").Value = 0.0009631928586580701
$ws.Range("S11").Value = 0.255338915736973
$ws.Range("T11").Value = 0.0001387130736508046
$ws.Range("U11").Value = 0.2564753741715639
$ws.Range("B12").Value = -55.96025251709802
$ws.Range("C12").Value = 0.001278282583173604
$ws.Range("H12").Value = 2.733406780273269
$ws.Range("I12").Value = 0.0309867296100948
$ws.Range("L12").Value = 0.007155010179999642
$ws.Range("M12").Value = 0.1354055892860337
$ws.Range("N12").Value = 0.0000518926478630097
$ws.Range("O12").Value = 0.1354055892860333
$ws.Range("P12").Value = 1.318059551970325
$ws.Range("Q12").Value = 0.4537665522382115
$ws.Range("T12").Value = 0.2343682381076229
$ws.Range("U12").Value = 1.283193429135786
